# Applies the "Updated symbol list" crypto-price refresh described in the commit.
# Columns D (Price) and E (Volume(1h)) hold numeric/percent-looking text that must
# stay plain text (as it was authored), so those assignments are apostrophe-prefixed
# to force Excel to keep them as text instead of re-typing them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''259.31'
$ws.Range('E2').Value = '''0.77%'

$ws.Range('D3').Value = '''27.03'
$ws.Range('E3').Value = '''0.21%'

$ws.Range('D4').Value = '''4.708'
$ws.Range('E4').Value = '''0.55%'

$ws.Range('D5').Value = '''0.06031'
$ws.Range('E5').Value = '''2.48%'

$ws.Range('D6').Value = '''6.675'
$ws.Range('E6').Value = '''0.42%'

$ws.Range('D7').Value = '''0.8599'
$ws.Range('E7').Value = '''0.27%'

$ws.Range('D8').Value = '''0.9243'
$ws.Range('E8').Value = '''-3.18%'

$ws.Range('E9').Value = '''-0.63%'

$ws.Range('D10').Value = '''0.05061'
$ws.Range('E10').Value = '''28.67%'

$ws.Range('D11').Value = '''0.07079'
$ws.Range('E11').Value = '''-0.22%'

$ws.Range('D12').Value = '''0.03084'
$ws.Range('E12').Value = '''-3.11%'

$ws.Range('D13').Value = '''0.09130'
$ws.Range('E13').Value = '''-0.43%'

$ws.Range('D14').Value = '''0.001530'
$ws.Range('E14').Value = '''-1.22%'

$ws.Range('D15').Value = '''0.0006076'
$ws.Range('E15').Value = '''0.75%'

$ws.Range('D16').Value = '''0.006054'
$ws.Range('E16').Value = '''-2.42%'

$ws.Range('E17').Value = '''-1.43%'

$ws.Range('D18').Value = '''3.170'

$ws.Range('D19').Value = '''2.165'
$ws.Range('E19').Value = '''-2.78%'

$ws.Range('E20').Value = '''0.41%'

$ws.Range('D21').Value = '''0.1298'
$ws.Range('E21').Value = '''0.42%'

$ws.Range('D22').Value = '''4.123'
$ws.Range('E22').Value = '''7.01%'

$ws.Range('D23').Value = '''0.04239'
$ws.Range('E23').Value = '''0.03%'

$ws.Range('E24').Value = '''-0.49%'

$ws.Range('D25').Value = '''0.004036'
$ws.Range('E25').Value = '''-6.11%'

$ws.Range('D26').Value = '''0.0001200'
$ws.Range('E26').Value = '''-0.04%'

$ws.Range('D27').Value = '''0.0001524'
$ws.Range('E27').Value = '''-21.35%'

$ws.Range('D40').Value = '''0.03850'
$ws.Range('E40').Value = '''0.52%'

$ws.Range('E41').Value = '''0.99%'

$ws.Range('D42').Value = '''0.004016'
$ws.Range('E42').Value = '''-36.17%'

$ws.Range('D43').Value = '''0.01527'
$ws.Range('E43').Value = '''33.71%'

$ws.Range('D44').Value = '''0.002200'
$ws.Range('E44').Value = '''-9.84%'

$ws.Range('D45').Value = '''0.00005115'
$ws.Range('E45').Value = '''-5.99%'

$ws.Range('E46').Value = '''0.00%'

$ws.Range('B47').Value = 'BOLO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range('D47').Value = '''0.1321'
$ws.Range('E47').Value = '''-21.15%'

$ws.Range('B48').Value = 'CoinbaseStockToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range('D48').Value = '''0.05456'
$ws.Range('E48').Value = '''-9.08%'

$ws.Range('D49').Value = '''0.00002100'
$ws.Range('E49').Value = '''0.00%'

$ws.Range('D50').Value = '''0.0002000'
$ws.Range('E50').Value = '''0.00%'
